$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11175
$ws.Range("C2").Value = 550764.5

$ws.Range("B3").Value = 116155.8
$ws.Range("C3").Value = 1469916.49

$ws.Range("B4").Value = 533076.84
$ws.Range("C4").Value = 2255888.8

$ws.Range("B5").Value = 1397833.49
$ws.Range("C5").Value = 3452281.04

$ws.Range("B6").Value = 1625135.07
$ws.Range("C6").Value = 3499868.84

$ws.Range("B7").Value = 1087507.69
$ws.Range("C7").Value = 3382567.31

$ws.Range("B8").Value = 2702868.57
$ws.Range("C8").Value = 3800337.7

$ws.Range("B9").Value = 1299503.9
$ws.Range("C9").Value = 1717223
